# Appends two new paragraphs at the end of the document body:
#   1) a highlighted (yellow) Times New Roman paragraph containing three
#      leading tabs followed by the Latin text
#      "Eamus omnes. Vos valete, spectatores optimi,"
#   2) an empty "Zkladntext" (Body Text) styled paragraph.
#
# Built as a raw OOXML fragment and inserted via Range.InsertXML so the
# resulting w:p/w:r markup (run splits, rFonts/sz/szCs/highlight on both
# the runs and the paragraph mark) matches exactly.

$d = $word.ActiveDocument

$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:tab/>
        <w:t xml:space="preserve">Eamus </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>o</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>mnes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>V</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:sz w:val="24"/>
          <w:szCs w:val="24"/>
          <w:highlight w:val="yellow"/>
        </w:rPr>
        <w:t>os valete, spectatores optimi,</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Zkladntext"/>
      </w:pPr>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint = $d.Range($d.Content.End, $d.Content.End)
$insertionPoint.InsertXML($xmlFrag)

Write-Output ("paragraphs now: " + $d.Paragraphs.Count)
